$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Output")
$ws.Range("J5").Value = 2772.237396664884
$ws.Range("K5").Value = 3495.002032033671
$ws.Range("F9").Value = 0.002999999997882697
$ws.Range("G9").Value = 194.5056000000021
$ws.Range("G12").Value = 3106.526361903447
$ws.Range("H12").Value = 229.3026380965508
$ws.Range("G13").Value = 1816.240638096552
$ws.Range("J13").Value = 385.2573619034488
$ws.Range("H15").Value = 442.7893619034492
$ws.Range("I18").Value = 1300.803101502453
$ws.Range("J18").Value = 1198.480144735196
$ws.Range("K18").Value = 394.3357277110096
$ws.Range("L18").Value = 346.1032977152996
$ws.Range("D20").Value = 10573.37281062786
$ws.Range("E20").Value = 9813.491920798784
$ws.Range("F20").Value = 8330.126451514641
$ws.Range("G20").Value = 8877.49174572686
$ws.Range("H20").Value = 8365.601310035541
$ws.Range("I20").Value = 7853.71087434422
$ws.Range("K20").Value = 6673.549543797599
$ws.Range("D21").Value = 467.8900000000004
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 16.27510747906919
$ws.Range("K21").Value = 11.67210208608054
$ws.Range("D22").Value = 373.6
$ws.Range("E22").Value = 521.9999999999999
$ws.Range("G22").Value = 809.8999999999999
$ws.Range("H22").Value = 949.3999999999996
$ws.Range("I22").Value = 972.6000000000004
$ws.Range("J22").Value = 995.8000000000001
$ws.Range("D23").Value = 110.96307332494
$ws.Range("E23").Value = 138.2166017570551
$ws.Range("F23").Value = 198.6000000000001
$ws.Range("I23").Value = 288.3000000000001
$ws.Range("J23").Value = 288.4999999999999
$ws.Range("L23").Value = 284.8
$ws.Range("D24").Value = 346.6000000000001
$ws.Range("F24").Value = 603.2
$ws.Range("G24").Value = 742.8499999999999
$ws.Range("I24").Value = 913.5000000000001
$ws.Range("K24").Value = 969.8000000000004
$ws.Range("L24").Value = 995.1000000000001
$ws.Range("E25").Value = 272.7276093977164
$ws.Range("L25").Value = 1667.7
$ws.Range("G27").Value = 93.95
$ws.Range("F28").Value = 668.3862529495993
$ws.Range("D29").Value = 116.0285566750598
$ws.Range("E29").Value = 76.34841884522939
$ws.Range("F32").Value = 933.6
$ws.Range("G32").Value = 940.2499999999999
$ws.Range("I32").Value = 945.3999999999996
$ws.Range("D33").Value = 424.3000000000001
$ws.Range("E33").Value = 447.65
$ws.Range("F33").Value = 471
$ws.Range("H33").Value = 565
$ws.Range("I33").Value = 610.15
$ws.Range("J33").Value = 655.2999999999998
$ws.Range("K33").Value = 707.0999999999999
$ws.Range("L33").Value = 758.8999999999999
$ws.Range("G34").Value = 31.35
$ws.Range("H34").Value = 34.9
$ws.Range("J34").Value = 43
$ws.Range("K34").Value = 48.3
$ws.Range("F35").Value = 1022.347377050401
$ws.Range("G36").Value = 71.40000000000002
$ws.Range("I36").Value = 65.00000000000001
$ws.Range("K36").Value = 56.39999999999998
$ws.Range("F37").Value = 270.0005
$ws.Range("G37").Value = 202.8997751480056
$ws.Range("H37").Value = 468.999299999997
$ws.Range("I37").Value = 0.03849999999925556
$ws.Range("K37").Value = 0.0385
$ws.Range("L37").Value = 0.03849999999999999
$ws.Range("G38").Value = 161.4110248519913
$ws.Range("H38").Value = 51.06128713508157
$ws.Range("I38").Value = 46.74522442266372
$ws.Range("E39").Value = 283
$ws.Range("F39").Value = 283
$ws.Range("G40").Value = 2273.395691447231
$ws.Range("H40").Value = 2323.176350186845
$ws.Range("I40").Value = 2390.119196497897
$ws.Range("J40").Value = 2470.436355501527
$ws.Range("K40").Value = 2557.089759884338
$ws.Range("L40").Value = 2708.536957775339
$ws.Range("D41").Value = 166.0537997527061
$ws.Range("E41").Value = 396.1279751793557
$ws.Range("F41").Value = 520.3857758013409
$ws.Range("G41").Value = 630.3037656833579
$ws.Range("H41").Value = 736.2677765340173
$ws.Range("I41").Value = 857.3430349814155
$ws.Range("J41").Value = 1023.631040539374
$ws.Range("K41").Value = 1267.548074091485
$ws.Range("L41").Value = 1427.457391435527
$ws.Range("D42").Value = 164.23
$ws.Range("E42").Value = 361.6799999999999
$ws.Range("G42").Value = 1187.98602997674
$ws.Range("H42").Value = 1894.015241261687
$ws.Range("I42").Value = 2741.899210747633
$ws.Range("J42").Value = 3635.688880143752
$ws.Range("K42").Value = 4522.191697904803
$ws.Range("L42").Value = 5158.046184530204
$ws.Range("D43").Value = 325.0012002472943
$ws.Range("E43").Value = 465.3540248206446
$ws.Range("F43").Value = 631.9832241986604
$ws.Range("G43").Value = 818.7252309609955
$ws.Range("H43").Value = 1449.427064294329
$ws.Range("I43").Value = 2386.546001547121
$ws.Range("J43").Value = 4394.286666278605
$ws.Range("K43").Value = 5535.86488906396
$ws.Range("L43").Value = 6521.922795640871
$ws.Range("D44").Value = 2.74
$ws.Range("F44").Value = 20.45999999999867
$ws.Range("G44").Value = 17.5512894533074
$ws.Range("H44").Value = 18.31493590894711
$ws.Range("I44").Value = 26.25180673764283
$ws.Range("J44").Value = 56.82934348650203
$ws.Range("K44").Value = 57.40330690548965
$ws.Range("D45").Value = 0.01448684999999994
$ws.Range("E45").Value = 0.6680376000000001
$ws.Range("F45").Value = 1.2725375184
$ws.Range("G45").Value = 1.3897142568
$ws.Range("H45").Value = 1.363098924
$ws.Range("I45").Value = 0.7542386280000002
$ws.Range("J45").Value = 2.77863696
$ws.Range("K45").Value = 2.967852960000001
$ws.Range("L45").Value = 3.185136
$ws.Range("D46").Value = 1.70551315
$ws.Range("E46").Value = 1.0919624
$ws.Range("F46").Value = 0.5774624816000002
$ws.Range("D47").Value = 550.0799999999999
$ws.Range("E47").Value = 983.0199999999999
$ws.Range("G47").Value = 2467.489409752511
$ws.Range("H47").Value = 3058.358621344633
$ws.Range("I47").Value = 3253.540558178556
$ws.Range("J47").Value = 3411.542161175159
$ws.Range("K47").Value = 3530.689433333829
$ws.Range("L47").Value = 3542.75503808367
$ws.Range("F48").Value = 259.57
$ws.Range("G48").Value = 825.1520544
$ws.Range("H48").Value = 1390.652263331941
$ws.Range("I48").Value = 1995.982179335768
$ws.Range("J48").Value = 2397.2497435152
$ws.Range("K48").Value = 2443.730157504144
$ws.Range("L48").Value = 2455.6127596128
$ws.Range("D50").Value = 45.454
$ws.Range("L52").Value = 0
$ws.Range("G55").Value = 8951.877980473259
$ws.Range("L55").Value = 2497.620175446785
$ws.Range("H56").Value = 10387.68811051881
$ws.Range("I56").Value = 1626.787953546271
$ws.Range("D57").Value = 5.683000000000001
$ws.Range("G60").Value = 12494.9849916314
$ws.Range("H60").Value = 5058.719403321727
$ws.Range("I60").Value = 1593.302477607718
$ws.Range("K61").Value = 2283.527999999999
$ws.Range("L61").Value = 2728.880956032416
$ws.Range("G62").Value = 3833.038814999998
$ws.Range("G63").Value = 493.7868
$ws.Range("H63").Value = 703.4
$ws.Range("J63").Value = 36.62799999999994
$ws.Range("K63").Value = 36.62800000000001
$ws.Range("L63").Value = 427.8341240396194
$ws.Range("H64").Value = 386.8699999999998
$ws.Range("I64").Value = 466.0024999999998
$ws.Range("D65").Value = 8181.324
$ws.Range("F65").Value = 7355.055
$ws.Range("G65").Value = 4512.146192142857
$ws.Range("J65").Value = 855.2714999999999
$ws.Range("K65").Value = 400.3398510638298
$ws.Range("F67").Value = 978.6645
$ws.Range("D73").Value = 586.5879179999999
